{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n// \"\u00a9 2020 . Contact: ...\" copyright/footer line, and the blank paragraph that\n// separated them from the \"Requisitos\" section, while leaving the single\n// blank paragraph that originally followed them (right before the final\n// page-break paragraph) intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4205: ...\" requirement paragraph; the three paragraphs to\n// delete immediately follow it (blank, \"Ver no Jupiter ...\", \"\u00a9 2020 ...\").\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4205\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4205' requirements paragraph.\");\n}\n\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i < items.length; i++) {\n  const text = items[i].text;\n  if (\n    text === \"\" ||\n    text.indexOf(\"Ver no Jupiter\") !== -1 ||\n    text.indexOf(\"\\u00a9 2020\") !== -1\n  ) {\n    toDelete.push(items[i]);\n    // Stop once we've captured the blank separator + both footer lines.\n    if (text.indexOf(\"\\u00a9 2020\") !== -1) {\n      break;\n    }\n  } else {\n    break;\n  }\n}\n\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n# \"\u00a9 2020 . Contact: ...\" copyright/footer line, and the blank paragraph that\n# separated them from the \"Requisitos\" section, while leaving the single\n# blank paragraph that originally followed them (right before the final\n# page-break paragraph) intact.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n\n# Locate the \"LOQ4205: ...\" requirement paragraph; the three paragraphs to\n# remove immediately follow it (blank, \"Ver no Jupiter ...\", \"\u00a9 2020 ...\").\n$anchor = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*LOQ4205*\") {\n        $anchor = $i\n        break\n    }\n}\n\nif ($anchor -eq -1) {\n    throw \"Could not find the 'LOQ4205' requirements paragraph.\"\n}\n\n# Collect the run of paragraphs right after the anchor: an empty separator\n# paragraph, the \"Ver no Jupiter ...\" line and the \"\u00a9 2020 ...\" line. Stop as\n# soon as the footer line is consumed so an unexpected layout doesn't cause\n# extra paragraphs to be swept up.\n$last = -1\nfor ($i = $anchor + 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    $core = $t.TrimEnd(\"`r\", \"`a\")\n    if ($core -eq \"\" -or $core -like \"*Ver no Jupiter*\" -or $core -like \"*2020*Contact*\") {\n        $last = $i\n        if ($core -like \"*2020*Contact*\") {\n            break\n        }\n    } else {\n        break\n    }\n}\n\nif ($last -eq -1) {\n    throw \"Could not find the paragraphs to remove after 'LOQ4205'.\"\n}\n\n$first = $anchor + 1\n$startRange = $d.Paragraphs($first).Range\n$endRange = $d.Paragraphs($last).Range\n$deleteRange = $d.Range($startRange.Start, $endRange.End)\n$deleteRange.Delete()\n"}
